$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.235323905944824
$ws.Range("B1").Value = 1.425225257873535
$ws.Range("C1").Value = 1.181261897087097
$ws.Range("D1").Value = 1.198772788047791
$ws.Range("E1").Value = 1.123106122016907
